$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$t = $m.Theme
$cs = $t.ThemeColorScheme

# Swap the presentation's (Integral) theme colour scheme for the
# stock "Office" colour scheme -- i.e. the colours that currently
# live in ppt/theme/theme1.xml become the ones driving the slide
# master's theme part (ppt/theme/theme2.xml).
$cs.Item(1).RGB  = 0          ' dk1      000000
$cs.Item(2).RGB  = 16777215   ' lt1      FFFFFF
$cs.Item(3).RGB  = 6968388    ' dk2      44546A
$cs.Item(4).RGB  = 15132391   ' lt2      E7E6E6
$cs.Item(5).RGB  = 13998939   ' accent1  5B9BD5
$cs.Item(6).RGB  = 3243501    ' accent2  ED7D31
$cs.Item(7).RGB  = 10855845   ' accent3  A5A5A5
$cs.Item(8).RGB  = 49407      ' accent4  FFC000
$cs.Item(9).RGB  = 12874308   ' accent5  4472C4
$cs.Item(10).RGB = 4697456    ' accent6  70AD47
$cs.Item(11).RGB = 12673797   ' hlink    0563C1
$cs.Item(12).RGB = 7491477    ' folHlink 954F72
